$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Through 2022-04-14")

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-15"

# Update header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 04-15)"

# Update the May 2022 value (row 5, column I) and recompute the Total row (row 14)
$ws.Range("I5").Value = 62
$ws.Range("I14").Value = 496
